$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New project items (task name -> date serial) appended after the existing
# rows (rows 25-34), continuing the same table started in row 3.
$items = @(
    @{ Row = 25; Name = "Integração da API no tela de login e cadastro"; Date = 45427 },
    @{ Row = 26; Name = "Manual de Instalação "; Date = 45433 },
    @{ Row = 27; Name = "Arduino conectado no Banco de Dados"; Date = 45455 },
    @{ Row = 28; Name = "Site institucional - Versão Final"; Date = 45456 },
    @{ Row = 29; Name = "Dashboards Conectadas"; Date = 45456 },
    @{ Row = 30; Name = "Fluxogramas de processos de atendimento"; Date = 45432 },
    @{ Row = 31; Name = "Ferramenta de Help Desk"; Date = 45440 },
    @{ Row = 32; Name = "Documentação de GMUD"; Date = 45438 },
    @{ Row = 33; Name = "Modelagem do Banco de Dados - Versão Final"; Date = 45416 },
    @{ Row = 34; Name = "Distribuição dos Servidores locais em três máquinas"; Date = 45425 }
)

foreach ($item in $items) {
    $r = $item.Row

    # Task name in column D: copy formatting (plain border style, same as
    # the rest of the table) from an existing cell, then set the text.
    $ws.Range("D20").Copy()
    $ws.Cells.Item($r, 4).PasteSpecial(-4122)
    $ws.Cells.Item($r, 4).Value = $item.Name

    # Date in column E: copy the date number-format/border/font style from an
    # existing date cell in the table, then set the value.
    $ws.Range("E20").Copy()
    $ws.Cells.Item($r, 5).PasteSpecial(-4122)
    $ws.Cells.Item($r, 5).Value = $item.Date
}

# Two additional blank rows (35-36) matching the existing blank-row style.
$ws.Range("D20").Copy()
$ws.Range("D35:E36").PasteSpecial(-4122)

# Update the view: zoom out to 70% and move the active selection.
$ws.Application.ActiveWindow.Zoom = 70
$ws.Range("G37").Select()
